$d = $word.ActiveDocument

# =================================================================
# PART 1 — split the "Hoy en día..." paragraph: insert three new,
# empty paragraphs (inheriting the same "jc=both" formatting) right
# before it.
# =================================================================
for ($i = 0; $i -lt 3; $i++) {
    $rng = $d.Content
    $rng.Find.Execute("Hoy en día podemos encontrarla")
    $target = $rng.Paragraphs(1).Range
    $split = $target.Duplicate
    $split.Collapse(1)
    $split.Text = "`r"
}

# =================================================================
# PART 2 — grow the paragraph that used to hold only the _GoBack
# bookmark into the long new block of text (several new paragraphs).
# We append everything as one chunk (paragraph breaks encoded as
# "`r") right after the existing (empty) paragraph content, i.e.
# right after the bookmark tags, so nothing else has to move.
# =================================================================
$rng = $d.Content
$rng.Find.Execute("Los modelos gráficos probabilísticos no siempre")
$bmPara = $rng.Paragraphs(1).Range.Next(4)     # story-relative next paragraph (the bookmark one)
$bm = $d.Bookmarks("_GoBack")
$bmParaRange = $bm.Range.Paragraphs(1).Range
$insertPoint = $d.Range($bmParaRange.Start, $bmParaRange.End - 1)
$insertPoint.Collapse(0)

$block2 = ""
$block2 += "`tCon la tecnología van creciendo sus herramientas para desarrollarla, el desarrollo de algoritmos y circuitos especializados van mejorando día a día lo cual nos permite acceder a mejores formas de desarrollar. Uno de estos productos es la inteligencia artificial, la cual mejora día a día para ser implementada en muchas áreas, como lo es la ciencia de los datos, quien engloba a la Big Data. Cada día millones de datos son subidos a la nube."
$block2 += "`r`tEn una empresa se deben tomar en cuenta muchos factores a la hora de vender un productor u ofrecer un servicio. Al ser empresas grandes deben ser muy cautelosas a la hora de tomar decisiones. La AI puede evaluar factores y computar una respuesta eficiente para la empresa y tener seguras sus ventas y servicios. En cuanto a manufactura, hay robots que puedes desplazarse y localizar ciertas cosas. En la agricultura analizar las cosechas y el suelo, entre otros sectores."
$block2 += "`rAsí como la IA proveerá de empleos a personas que sepan manejar datos, programar y otros conocimientos de computación, eliminará la necesidad de contratar personas cuyo trabajo puede ser automatizable."
$block2 += "`r`tA todo esto respecto a la inteligencia artificial, se debe tener cuidado al momento de ser desarrollada, puesto que puede fallar, o simplemente que el creador no tome en cuenta ciertos criterios y conduzca a fracasar en la toma de alguna decisión a la máquina. Por ejemplo, anteriormente se ejemplificaba cómo una empresa podría tomar decisiones por medio de una IA, la cual podría fracasar y hacer perder a la empresa millones de dólares. Los sistemas que más están expuestos al peligro son los que toman datos de internet, pues pueden aprender prejuicios, comportamientos racistas, xenófobos entre otros."
$block2 += "`rConclusión."
$block2 += "`r`tLa inteligencia artificial ha crecido exponencialmente en los últimos años. De nacer en una enorme, monstruosa y tosca máquina como La Máquina de Turing, ahora lo podemos ver en pequeños chips que pueden caber donde sea, lo vemos en nuestros teléfonos, en nuestros relojes, en nuestras computadoras, por lo que podríamos definir la inteligencia artificial como omnipresente. Esta misma capacidad de estar en todos lados es lo que permite a la misma poder cuantificar, analizar y tomar decisiones en distintos sectores, en industrias y empresas grandes. Cada vez se necesitarán más personas capacitadas para manejar/manipular datos, y otros ingenieros en el campo de la computación. A su vez, se disminuirán costos de mano de obra automatizando estos puestos."
$block2 += "`r`tComo un niño, se le debe enseñar a la inteligencia artificial qué es bueno aprender y qué no, pues están en peligro de aprender conductas no morales las cuales pueden lastimar a minorías, y otros comportamientos."

$insertPoint.InsertAfter($block2)

Write-Output ("Paragraph count after part 2: " + $d.Paragraphs.Count)
